$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Summary section: trim the first paragraph and split off a new BodyText
#    paragraph with the remaining content.
# ---------------------------------------------------------------------------
$oldSummary = "Ryan Parman is a cloud-native engineering leader with over 25 years of experience, who specializes in technical leadership, software development, site reliability engineering, and cybersecurity for the modern web. Excels at listening, adapting, and driving continuous improvement. Delivers exceptional work, builds impactful solutions, and elevates team performance."
$newSummary1 = "Ryan Parman is a cloud-native engineering leader, who specializes in technical leadership, software development, site reliability engineering, and cybersecurity for the modern web. Excels at listening, adapting, and driving continuous improvement."
$newSummary2 = "Small business owner, two-time startup founder, and creator of two open-source projects with millions of users each. Ryan has a proven track record of building high-quality software, delivering impactful solutions, and elevating team performance."

$rng = $d.Content
$found = $rng.Find.Execute($oldSummary, $true, $false, $false, $false, $false, $true, 1, $false, ($newSummary1 + "^p" + $newSummary2), 2)

$rng2 = $d.Content
$found2 = $rng2.Find.Execute($newSummary2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2.Paragraphs(1).Style = "BodyText"

# ---------------------------------------------------------------------------
# 2) McGraw Hill: add a BlockText description paragraph after the
#    "McGraw Hill ... Remote (since COVID) ..." heading line.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("— Remote (since COVID), previously Seattle, WA", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.Move(4, 1) | Out-Null
$rng.Paragraphs(1).Style = "BlockText"

$rng.InsertAfter("McGraw Hill is a")
$rng.Collapse(0)
$rng.InsertAfter(" ")
$rng.Collapse(0)
$rng.InsertAfter("learning science")
$rng.Italic = $true
$rng.Collapse(0)
$rng.InsertAfter(" ")
$rng.Collapse(0)
$rng.InsertAfter(" company which produces textbooks, digital learning tools, and adaptive technology to enhance learning. It is one of the “big three” educational publishers in the U.S, and was acquired by Platinum Equity 2021.")

# ---------------------------------------------------------------------------
# 3) New bullet under "Principal Cloud and Platform Engineer" — inserted
#    as the first numbered item, ahead of "Managed the Base AMI ..." bullet.
# ---------------------------------------------------------------------------
$newBullet = "As every school in America transitioned to online learning during the COVID-19 lockdowns, I was the technical/development lead on the team who supported all SRE and product engineering teams, working on core platforms and services."
$rng = $d.Content
$found = $rng.Find.Execute("Managed the Base", $true, $false, $false, $false, $false, $true, 2, $false, ($newBullet + "^pManaged the Base"), 2)

# ---------------------------------------------------------------------------
# 4) WePay: add a BlockText description paragraph after the
#    "WePay ... Redwood City, CA" heading line.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("— Redwood City, CA", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.Move(4, 1) | Out-Null
$rng.Paragraphs(1).Style = "BlockText"
$rng.InsertAfter("WePay is an online payment service provider which provides “payments for platforms”, where examples of platforms are GoFundMe, Care.com, and Xbox. It was acquired by JPMorgan Chase in October 2017.")

# ---------------------------------------------------------------------------
# 5) Amazon Web Services: add a BlockText description paragraph after the
#    "Amazon Web Services ... Seattle, WA" heading line.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("— Seattle, WA", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.Move(4, 1) | Out-Null
$rng.Paragraphs(1).Style = "BlockText"
$rng.InsertAfter("Amazon Web Services provides on-demand cloud computing platforms and APIs to individuals, companies, and governments, on a metered, pay-as-you-go basis.")

Write-Host "All edits applied."
